$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gameflow")

# Update dialogue text in C10 and C11
$ws.Range("C10").Value = "Welcome back Uncle Nintendo! You last activated this <color=blue>TYPOCRYPHA</color> unit 1 year, 1 month, and 13 days ago. You have <color=yellow>1023</color> new updates and <color=yellow>33,333</color> new emails."
$ws.Range("C11").Value = "O-Oh. No, Im not him. I-uh . . . Dammit! How do I reset-"

# Update the selected/active cell to C10
$ws.Activate()
$ws.Range("C10").Select()
